$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = [double]"0.9999999667885955"
$ws.Range("E2").Value = [double]"0.9999999667885955"

$ws.Range("D3").Value = [double]"0.01462120494210556"
$ws.Range("E3").Value = [double]"0.01462120494210556"

$ws.Range("D4").Value = [double]"4.352120057555408E-10"
$ws.Range("E4").Value = [double]"4.352120057555408E-10"

$ws.Range("D5").Value = [double]"1.691215890752005E-29"
$ws.Range("E5").Value = [double]"1.691215890752005E-29"

$ws.Range("D6").Value = [double]"3.905607788615773E-50"
$ws.Range("E6").Value = [double]"3.905607788615773E-50"

$ws.Range("D8").Value = [double]"0.999999999997889"
$ws.Range("E8").Value = [double]"2.110978059022273E-12"

$ws.Range("D9").Value = 0

$ws.Range("D10").Value = [double]"5.492476505148903E-06"
$ws.Range("E10").Value = [double]"0.9999945075234948"

$ws.Range("D11").Value = [double]"0.0003789383652183232"
$ws.Range("E11").Value = [double]"0.9996210616347817"
$ws.Range("F11").Value = [double]"99.01937103271484"
